# "arreglos de la parte de excel"
# Update the four column-header cells in row 5 of Hoja1 with the expanded
# wording, which grows row 5's (wrap-text) height, and move the current
# selection from H5 to G6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)
$ws.Activate()

$ws.Range("C5").Value = "Raised in`n(Soft. Vers.)"
$ws.Range("E5").Value = "Fixed in`n(Soft. Vers.)"
$ws.Range("F5").Value = "Fixed by QA or on prototype (prototype by default)"
$ws.Range("G5").Value = "Tested by  MCO"

# Row 5 grows a little taller to fit the new wrapped text (67.8 -> 69.6 pt).
$ws.Rows.Item(5).RowHeight = 69.6

# Selection moves from H5 to G6, and the view scrolls so row 3 is at the top.
$ws.Range("G6").Select()
